$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing three data rows (originally rows 20-22) first; the
# runtime mishandles a delete that follows an insert, so the delete must
# happen before the insert below.
$ws.Range("A20:A22").EntireRow.Delete()

# Insert two new rows right after the header row, shifting the remaining
# data rows down to make room for the new samples.
$ws.Range("A2:A3").EntireRow.Insert()

# The inserted rows pick up formatting from their neighbours; clear it so
# the new data rows stay unstyled like every other data row.
$ws.Range("A2:C3").ClearFormats()

# Populate the two newly inserted rows with their values.
$ws.Range("A2").Value = -0.0704022198915481
$ws.Range("B2").Value = 0.1944078654050827
$ws.Range("C2").Value = -0.0245873257517814

$ws.Range("A3").Value = 0.3005456924438476
$ws.Range("B3").Value = 0.8894197940826416
$ws.Range("C3").Value = 0.086895577609539
